# goa & mp done
#
# - files_info: add a new "goa" row (folder_location / filename / date_format),
#   mirroring the existing gujarat / mp rows.
# - meas_info: point the "goa" demand_met address at AD8,AD9,AD10,AD11
#   (was Z8,Z9,Z10,Z11) and the "mp" demand_met_mu address at E41 (was A1).
# - Leave the workbook focused on meas_info (scrolled near the bottom,
#   row 51 selected) instead of files_info.

$wb = $excel.ActiveWorkbook

# --- files_info: add a new "goa" row ---
$filesInfo = $wb.Worksheets.Item("files_info")
$filesInfo.Range("A4").Value = "goa"
$filesInfo.Range("B4").Value = "C:\Users\dheer\Desktop\wrldc\wrldc_mis_state_files_ingestion\stateFiles"
$filesInfo.Range("C4").Value = "GOA_{{dt}}_uploaded.xlsx"
$filesInfo.Range("D4").Value = "%d_%m_%Y"

# --- meas_info: fix up the goa / mp address cells ---
$measInfo = $wb.Worksheets.Item("meas_info")
$measInfo.Range("E50").Value = "AD8,AD9,AD10,AD11"
$measInfo.Range("E51").Value = "E41"

# --- view state: files_info keeps its old selection but is no longer the
#     active tab; meas_info becomes active, scrolled down, with E51 selected ---
$filesInfo.Range("E3").Select()

$measInfo.Activate()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$measInfo.Range("E51").Select()
